$d = $word.ActiveDocument
$d.Content.Find.Execute("Vignesh Sundaramahalingam", $true, $false, $false, $false, $false,
                         $true, 1, $false, " Vignesh Sundaramahalingam", 2)

# do an unrelated operation to clear any find-state
$junk = $d.Range(50,55)
Write-Output ("junk=" + $junk.Text)

$spaceRange = $d.Range(0,1)
$spaceRange.LanguageID = "en-US"
Write-Output "done"
